# chore: adapt column header formatting to respective input file names
#
# - Rename the "_old" / "_new" header-name suffixes to the respective
#   format-version suffixes "_FV2210" / "_FV2304".
# - Turn the data range into an actual Excel Table ("Table1").
# - Freeze the header row (row 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename headers: "<name>_old" -> "<name>_FV2210", "<name>_new" -> "<name>_FV2304" ---
$headerCols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U")
foreach ($col in $headerCols) {
    $cell = $ws.Range($col + "1")
    $text = [string]$cell.Value2
    if ($text.EndsWith("_old")) {
        $cell.Value = $text.Substring(0, $text.Length - 4) + "_FV2210"
    } elseif ($text.EndsWith("_new")) {
        $cell.Value = $text.Substring(0, $text.Length - 4) + "_FV2304"
    }
}

# --- 2. Turn A1:U57 into an Excel Table named "Table1" ---
# The header row already carries bold/centered/wrapped/filled/bordered
# formatting. Snapshot it into a scratch row first, strip the header
# row's own formatting so the engine doesn't bake it into a brand-new
# header-row dxf when the table is created, then paste the formatting
# back from the scratch row (re-using the existing style) and wipe the
# scratch row again.
$headerRange = $ws.Range("A1:U1")
$scratch = $ws.Range("A100:U100")

$headerRange.Copy()
$scratch.PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$headerRange.ClearFormats()

$tableRange = $ws.Range("A1:U57")
$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$tbl.Name = "Table1"

$scratch.Copy()
$headerRange.PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false
$scratch.Clear()

# --- 3. Freeze the header row (row 1) ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
